# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape snapshot (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 2234
    3  = 99
    4  = 13476
    7  = 529
    9  = 1200
    10 = 1013
    11 = 13843
    12 = 14576
    13 = 43
    21 = 50
    23 = 1123
    24 = 117
    25 = 60
    26 = 5586
    27 = 941
    28 = 1041
    29 = 5363
    30 = 38
    31 = 27
    32 = 172
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# ---- Sheet "全部类型" ----
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 2234
    3  = 99
    4  = 13476
    8  = 529
    10 = 1200
    11 = 1013
    12 = 13843
    13 = 14576
    14 = 43
    22 = 50
    24 = 1123
    25 = 117
    26 = 60
    27 = 5586
    28 = 941
    29 = 1041
    30 = 5363
    31 = 38
    32 = 27
    33 = 172
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
